$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 596, shifting all data
# currently in rows 596:618 down to rows 598:620.
$ws.Rows("596:597").Insert()

# --- New row 596 ---
$ws.Range("A596").Value = 11
$ws.Range("B596").Value = "Vega Monumental Concepción"
$ws.Range("C596").Value = "Bíobío"
$ws.Range("D596").Value = 45008
$ws.Range("E596").Value = 8
$ws.Range("F596").Value = "Fruta"
$ws.Range("G596").Value = 100104
$ws.Range("H596").Value = "Frutos de pepita"
$ws.Range("I596").Value = 100104005
$ws.Range("J596").Value = "Pera"
$ws.Range("K596").Value = "Abate Fettel"
$ws.Range("L596").Value = "Primera"
$ws.Range("M596").Value = 350
$ws.Range("N596").Value = 8000
$ws.Range("O596").Value = 9000
$ws.Range("P596").Value = 8571
$ws.Range("Q596").Value = "$/caja 16 kilos empedrada"
$ws.Range("R596").Value = "Provincia de Curicó"
$ws.Range("S596").Value = 536
$ws.Range("T596").Value = 16

# --- New row 597 ---
$ws.Range("A597").Value = 11
$ws.Range("B597").Value = "Vega Monumental Concepción"
$ws.Range("C597").Value = "Bíobío"
$ws.Range("D597").Value = 45008
$ws.Range("E597").Value = 8
$ws.Range("F597").Value = "Fruta"
$ws.Range("G597").Value = 100104
$ws.Range("H597").Value = "Frutos de pepita"
$ws.Range("I597").Value = 100104005
$ws.Range("J597").Value = "Pera"
$ws.Range("K597").Value = "Packham's Triumph"
$ws.Range("L597").Value = "Primera"
$ws.Range("M597").Value = 270
$ws.Range("N597").Value = 10000
$ws.Range("O597").Value = 11000
$ws.Range("P597").Value = 10556
$ws.Range("Q597").Value = "$/caja 16 kilos empedrada"
$ws.Range("R597").Value = "Provincia de Curicó"
$ws.Range("S597").Value = 660
$ws.Range("T597").Value = 16
